$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.004.03"
$ws.Range("E2").Value = "  +0.56%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.641.48"
$ws.Range("E3").Value = "  +0.67%  "

$ws.Range("E4").Value = "  +0.36%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.87"

$ws.Range("E6").Value = "  +0.40%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.35%  "

$ws.Range("E8").Value = "  +0.29%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0639"
$ws.Range("E9").Value = "  +1.03%  "

$ws.Range("E10").Value = "  +0.26%  "

$ws.Range("E11").Value = "  +0.67%  "

$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.869.94"
$ws.Range("E12").Value = "  +0.77%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.674.82"
$ws.Range("E13").Value = "  +0.12%  "

$ws.Range("E14").Value = "  +0.28%  "

$ws.Range("E15").Value = "  -0.27%  "

$ws.Range("E16").Value = "  +1.08%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.35"
$ws.Range("E17").Value = "  +1.23%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.096.40"
$ws.Range("E18").Value = "  +0.87%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.00"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "194.78"
$ws.Range("E20").Value = "  +0.84%  "

$ws.Range("E21").Value = "  -0.65%  "

$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.20"
$ws.Range("E23").Value = "  -0.81%  "

$ws.Range("B24").Value = "Stellar"
$ws.Range("C24").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.131"
$ws.Range("E24").Value = "  +3.78%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.79"
$ws.Range("E25").Value = "  -1.82%  "

$ws.Range("E26").Value = "  +0.41%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "143.10"
$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("E28").Value = "  +0.55%  "

$ws.Range("E29").Value = "  +0.56%  "

$ws.Range("E30").Value = "  +0.93%  "

$ws.Range("E31").Value = "  -0.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.30"

$ws.Range("E33").Value = "  +1.30%  "

$ws.Range("E34").Value = "  -1.81%  "

$ws.Range("E35").Value = "  +1.77%  "

$ws.Range("E36").Value = "  +0.52%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.128.06"
$ws.Range("E37").Value = "  -0.75%  "

$ws.Range("E38").Value = "  -1.33%  "

$ws.Range("E39").Value = "  -0.18%  "

$ws.Range("E40").Value = "  +0.41%  "

$ws.Range("E41").Value = "  +0.49%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.32"
$ws.Range("E42").Value = "  +0.26%  "

$ws.Range("E43").Value = "  -0.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.779.30"
$ws.Range("E44").Value = "  +0.82%  "

$ws.Range("E45").Value = "  +4.49%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.65"
$ws.Range("E46").Value = "  +0.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0523"
$ws.Range("E47").Value = "  -0.37%  "

$ws.Range("E48").Value = "  +1.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.72"
$ws.Range("E49").Value = "  +1.53%  "

$ws.Range("E50").Value = "  -0.20%  "

$ws.Range("E51").Value = "  -0.48%  "
